# Recalibration update: refresh numeric outputs on Linear and NonLinear
# sheets, and the corresponding "abs_epsi_autocorr" / "sig2_1" array strings.

$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")

$wsLinear.Range("B2").Value = -4.539101350276637
$wsLinear.Range("B3").Value = 0.03349084872856418
$wsLinear.Range("B4").Value = 735690.0479906698
$wsLinear.Range("B5").Value = "[1.0, 0.16567847909221484, 0.1281243159000581, 0.14938375448098118, 0.14185896035401052, 0.11676657279734318, 0.16910623000292685, 0.2991519261102621, 0.17334870233848512, 0.12153998399684587, 0.12808171312419367, 0.16003076500856628, 0.12551361900991267, 0.2006996338757427, 0.2538695051065751, 0.180809279457473, 0.12390492236174899, 0.16127477530651582, 0.13027199286800537, 0.1432814466719841]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")

$wsNonLinear.Range("B4").Value = -22.09673263352916
$wsNonLinear.Range("B5").Value = 0.04069433339166005
$wsNonLinear.Range("B6").Value = 1167782.585152757
$wsNonLinear.Range("B7").Value = 23.16739802285751
$wsNonLinear.Range("B8").Value = -0.06042499595970936
$wsNonLinear.Range("B9").Value = 345794.175302273
$wsNonLinear.Range("B10").Value = "[1.0, 0.16274706957587315, 0.12622235400711498, 0.14670709094378107, 0.13924730503912697, 0.113786036845273, 0.16671364028750796, 0.29511476800432596, 0.1701464032170918, 0.11901367183959839, 0.12513737684978443, 0.15705002772575022, 0.1232365854725888, 0.19743090271976882, 0.24964186387266674, 0.17710888463840793, 0.12078280634898264, 0.15848743469033447, 0.12693291277772342, 0.14017452334956537]"
